# Workbook tracks a running history of bank-SMS / notification entries per
# month on the "2024" sheet, newest entry first. A new September entry was
# captured, so insert a fresh row above the existing September 9th entries
# (row 35) which pushes all the rows below it down by one, then populate
# the new row's September_Details / September_Date columns (R / S) with the
# newly captured entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new blank row at row 35 — everything from row 35 down shifts to
# row 36 down (matches dimension growing from A1:Y98 to A1:Y99).
$ws.Rows.Item(35).Insert()

# Populate the new row with the latest captured entry.
$ws.Range("R35").Value = "corporate internet share"
$ws.Range("S35").Value = "2024-09-09 11:10:39"
